$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the feature attribution: move "Nadjia" from B20 up to B14
# (she is now credited for the row-14 feature instead of the row-20 one).
$ws.Range("B14").Value2 = $ws.Range("B20").Value2
$ws.Range("B20").ClearContents()

# Reset the active selection back to the top-left cell.
$ws.Range("A1").Select()
